$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0.001
$ws.Range("K12").Value = 470
$ws.Range("L12").Value = 0.00235
